$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2020-10", 18, 52, 34.61538461538461),
    @("2020-11", 17, 45, 37.77777777777778),
    @("2020-12", 13, 45, 28.88888888888889),
    @("2021-01", 20, 36, 55.55555555555556),
    @("2021-02", 25, 95, 26.31578947368421),
    @("2021-03", 30, 77, 38.96103896103897),
    @("2021-04", 37, 94, 39.36170212765958),
    @("2021-05", 43, 89, 48.31460674157304),
    @("2021-06", 52, 91, 57.14285714285714),
    @("2021-07", 63, 122, 51.63934426229508),
    @("2021-08", 61, 100, 61),
    @("2021-09", 64, 125, 51.2),
    @("2021-10", 70, 128, 54.6875),
    @("2021-11", 79, 139, 56.83453237410072),
    @("2021-12", 78, 127, 61.41732283464567),
    @("2022-01", 79, 117, 67.52136752136752),
    @("2022-02", 85, 143, 59.44055944055944),
    @("2022-03", 91, 143, 63.63636363636363),
    @("2022-04", 89, 146, 60.95890410958904),
    @("2022-05", 77, 139, 55.39568345323741),
    @("2022-06", 88, 112, 78.57142857142857),
    @("2022-07", 81, 117, 69.23076923076923),
    @("2022-08", 71, 104, 68.26923076923077),
    @("2022-09", 76, 108, 70.37037037037037),
    @("2022-10", 69, 127, 54.33070866141733),
    @("2022-11", 74, 98, 75.51020408163265),
    @("2022-12", 75, 123, 60.97560975609756),
    @("2023-01", 75, 104, 72.11538461538461),
    @("2023-02", 78, 124, 62.90322580645162),
    @("2023-03", 80, 121, 66.11570247933885),
    @("2023-04", 92, 130, 70.76923076923077),
    @("2023-05", 107, 141, 75.88652482269504),
    @("2023-06", 103, 143, 72.02797202797203),
    @("2023-07", 101, 138, 73.18840579710145),
    @("2023-08", 108, 143, 75.52447552447552),
    @("2023-09", 115, 160, 71.875),
    @("2023-10", 115, 156, 73.71794871794873),
    @("2023-11", 110, 164, 67.07317073170732),
    @("2023-12", 99, 150, 66),
    @("2024-01", 105, 143, 73.42657342657343),
    @("2024-02", 124, 182, 68.13186813186813),
    @("2024-03", 121, 162, 74.69135802469135),
    @("2024-04", 136, 189, 71.95767195767195),
    @("2024-05", 143, 197, 72.58883248730965),
    @("2024-06", 141, 209, 67.46411483253588),
    @("2024-07", 151, 215, 70.23255813953489),
    @("2024-08", 166, 226, 73.45132743362832),
    @("2024-09", 172, 226, 76.10619469026548),
    @("2024-10", 169, 227, 74.44933920704845),
    @("2024-11", 151, 242, 62.39669421487604),
    @("2024-12", 123, 197, 62.43654822335025),
    @("2025-01", 145, 195, 74.35897435897436),
    @("2025-02", 159, 216, 73.61111111111111),
    @("2025-03", 142, 229, 62.00873362445415),
    @("2025-04", 101, 204, 49.50980392156863)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $item = $data[$i]
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
}
